# A new weekly price-report row was inserted into the "Fruta, Feria Lagunitas
# de Puerto Montt - Mandarina" sheet at row 454, pushing every subsequent row
# (old 454..565) down by one (new 455..566). This mirrors how the upstream
# "logica_diaria" consolidation appends each new observation at the top of
# its date-sorted block instead of strictly at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 454, shifting rows 454:565 down to 455:566.
$ws.Rows("454:454").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A454").Value = 4
$ws.Range("B454").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C454").Value = "Los Lagos"
$ws.Range("D454").Value = 45275
$ws.Range("E454").Value = 10
$ws.Range("F454").Value = "Fruta"
$ws.Range("G454").Value = 100102
$ws.Range("H454").Value = "Cítricos"
$ws.Range("I454").Value = 100102004
$ws.Range("J454").Value = "Mandarina"
$ws.Range("K454").Value = "Murcott"
$ws.Range("L454").Value = "Primera"
$ws.Range("M454").Value = 300
$ws.Range("N454").Value = 16000
$ws.Range("O454").Value = 16000
$ws.Range("P454").Value = 16000
$ws.Range("Q454").Value = "$/bandeja 10 kilos"
$ws.Range("R454").Value = "Región de O'Higgins"
$ws.Range("S454").Value = 1600
$ws.Range("T454").Value = 10
